# connections.xlsx: "table of protocols, table of keywords"
#
# - "keywords" sheet is rebuilt into a proper table: Keyword | Is Protocol |
#   Is Abstract | Meanings (plus a small secondary Keyword/Meanings table for
#   ZMQ/RabbitMQ over in F:G).
# - The old "Keyword to scan for definition" / "Techniques" scratch columns
#   that used to live off to the side of "keywords" move out into a new
#   worksheet named "Others".
# - "keywords" becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Others" worksheet, positioned after "keywords".
# ---------------------------------------------------------------------
$keywordsSheet = $wb.Worksheets("keywords")
$others = $wb.Worksheets.Add($null, $keywordsSheet)
$others.Name = "Others"

# ---------------------------------------------------------------------
# 2. Rebuild the "keywords" sheet as the protocol table.
# ---------------------------------------------------------------------
$kw = $wb.Worksheets("keywords")
$kw.Cells.Clear()

# Header row.
$kw.Range("A1").Value = "Keyword"
$kw.Range("D1").Value = "Meanings"
$kw.Range("B1").Value = "Is Protocol"
$kw.Range("A2").Value = "MQTT"
$kw.Range("D2").Value = "Message Queuing Telemetry Transport"
$kw.Range("B2").Value = "YES"
$kw.Range("A3").Value = "OPC UA"
$kw.Range("D3").Value = "Open Platform Communication Unified Architecture"
$kw.Range("C1").Value = "Is Abstract"
$kw.Range("C3").Value = "YES"
$kw.Range("A4").Value = "DTLS"
$kw.Range("D4").Value = "Datagram Transport Layer Security"
$kw.Range("B4").Value = "YES"
$kw.Range("A5").Value = "PPMP"
$kw.Range("D5").Value = "Production Performance Management Protocol"
$kw.Range("B5").Value = "YES"
$kw.Range("A6").Value = "ABAC"
$kw.Range("D6").Value = "Attribute-based access control"
$kw.Range("C6").Value = "YES"
$kw.Range("A7").Value = "APPSTACLE"
$kw.Range("D7").Value = "open standard APplication Platform for carS and TrAnsportation vehiCLEs"
$kw.Range("C7").Value = "YES"
$kw.Range("A8").Value = "AGL"
$kw.Range("D8").Value = "Automotive Grade Linux"
$kw.Range("A9").Value = "DDS"
$kw.Range("D9").Value = "Data Distribution Service"
$kw.Range("C9").Value = "YES"
$kw.Range("A10").Value = "WoT"
$kw.Range("D10").Value = "Web of Thing"
$kw.Range("B10").Value = "YES"
$kw.Range("C10").Value = "YES"
$kw.Range("A11").Value = "Calypso"
$kw.Range("D11").Value = "the standard for contactless tickets and cards"
$kw.Range("B11").Value = "YES"

# Secondary little Keyword/Meanings table (F:G) - message brokers.
$kw.Range("F1").Value = "Keyword"
$kw.Range("G1").Value = "Meanings"
$kw.Range("F2").Value = "ZMQ"
$kw.Range("G2").Value = "ZeroMQ"
$kw.Range("F3").Value = "RabbitMQ"
$kw.Range("G3").Value = "open-source message-broker software"

# Column widths for the new layout.
$kw.Columns("A").ColumnWidth = 14.63
$kw.Columns("B").ColumnWidth = 10.67
$kw.Columns("C").ColumnWidth = 10.34
$kw.Columns("D").ColumnWidth = 66.1
$kw.Columns("G").ColumnWidth = 43.05

# ---------------------------------------------------------------------
# 3. Fill in the "Others" worksheet with the content that used to sit in
#    the side columns of "keywords" (keyword-scan helper list + NLP
#    technique notes).
# ---------------------------------------------------------------------
$others.Range("A1").Value = "Keyword to scan for definition"
$others.Range("G1").Value = "Techniques"
$others.Range("A2").Value = "about"
$others.Range("G2").Value = "Latent Semantic Indexing"
$others.Range("A3").Value = "introduction"
$others.Range("G3").Value = "WordNet"
$others.Range("A4").Value = "proposal"
$others.Range("A5").Value = "concept"
$others.Range("A6").Value = "overview"
$others.Range("A7").Value = "What is"

# ---------------------------------------------------------------------
# 4. Selection / active-tab bookkeeping: "keywords" becomes the active
#    sheet (was "connection data"), with "Others" left on its default
#    selection.
# ---------------------------------------------------------------------
$others.Range("D13").Select()
$kw.Range("G9").Select()
$kw.Activate()
